# "implemented new partition method - link process method"
#
# The eGRID Hg emission-factor table (Sheet1) is rebuilt with values from the
# new partition/link-process method:
#   - Rows whose B:D values were all 0 under the old method lose their B:D
#     data entirely (row 2, which had no category label, disappears; rows
#     3/4/6/9/10/13/15/16 keep only their category label in column A).
#   - Rows that still have non-zero contributions get their solid/liquid/gas
#     mg_mwh figures (columns B/C/D) replaced with the newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows that become label-only (B:D cleared) -----------------------------
# Row 2 has no category label at all, so clearing B:D leaves it fully empty
# and it drops out of the sheet, just like the source diff.
$ws.Range("B2:D2").ClearContents() | Out-Null
$ws.Range("B3:D3").ClearContents() | Out-Null
$ws.Range("B4:D4").ClearContents() | Out-Null
$ws.Range("B6:D6").ClearContents() | Out-Null
$ws.Range("B9:D9").ClearContents() | Out-Null
$ws.Range("B10:D10").ClearContents() | Out-Null
$ws.Range("B13:D13").ClearContents() | Out-Null
$ws.Range("B15:D15").ClearContents() | Out-Null
$ws.Range("B16:D16").ClearContents() | Out-Null

# --- Rows with updated solid/liquid/gas mg_mwh values -----------------------
$ws.Range("B5").Value2 = 10.030263468804595
$ws.Range("C5").Value2 = 0.27466506468859281
$ws.Range("D5").Value2 = 8.9437722055330475

$ws.Range("B7").Value2 = 63.338268618236683
$ws.Range("C7").Value2 = 0.70615137532559336
$ws.Range("D7").Value2 = 25.282344481046188

$ws.Range("B8").Value2 = 25.850112398748671
$ws.Range("C8").Value2 = 1.734005357129687
$ws.Range("D8").Value2 = 11.830814715060846

$ws.Range("B11").Value2 = 32.325938310834488
$ws.Range("D11").Value2 = 29.90704759777244

$ws.Range("B12").Value2 = 40.154462470993749
$ws.Range("C12").Value2 = 0
$ws.Range("D12").Value2 = 28.01502419513093

$ws.Range("B14").Value2 = 13.035302049642638
$ws.Range("C14").Value2 = 0.40885826421658927
$ws.Range("D14").Value2 = 14.416476200941689

$ws.Range("B17").Value2 = 40.171409052504941
$ws.Range("C17").Value2 = 1.6438284503855169
$ws.Range("D17").Value2 = 34.620596696822972

$ws.Range("B18").Value2 = 36.841194855346359
$ws.Range("C18").Value2 = 0.43714356403512616
$ws.Range("D18").Value2 = 9.1204043587328627

$ws.Range("B19").Value2 = 4.1229085237382082
$ws.Range("C19").Value2 = 0
$ws.Range("D19").Value2 = 53.210752066847533

$ws.Range("B20").Value2 = 37.531410130375463
$ws.Range("C20").Value2 = 2.8726640255485316
$ws.Range("D20").Value2 = 4.8618491656354355

$ws.Range("B21").Value2 = 39.896672744705995
$ws.Range("C21").Value2 = 0
$ws.Range("D21").Value2 = 22.530939550568363

$ws.Range("B22").Value2 = 17.251835739868511
$ws.Range("C22").Value2 = 0
$ws.Range("D22").Value2 = 44.383413618563992

$ws.Range("B23").Value2 = 25.07653965274055
$ws.Range("D23").Value2 = 50.473121540651618

$ws.Range("B24").Value2 = 45.889629975037529
$ws.Range("D24").Value2 = 17.305087623994527

$ws.Range("B25").Value2 = 48.532755881630855
$ws.Range("C25").Value2 = 0
$ws.Range("D25").Value2 = 12.10134012994699

$ws.Range("B26").Value2 = 28.738183274811526
$ws.Range("C26").Value2 = 0.80978737303684145
$ws.Range("D26").Value2 = 3.1185263198997348

$ws.Range("B27").Value2 = 31.09612521834983
$ws.Range("C27").Value2 = 1.6054482547722724
$ws.Range("D27").Value2 = 18.868140124036955

$ws.Range("B28").Value2 = 22.693288454845842
$ws.Range("C28").Value2 = 1.0603899780532324
$ws.Range("D28").Value2 = 9.4312119869575231

# --- Tidy up the active selection so it no longer points at the old C22 ----
# (the saved sheet view previously had the cursor parked on C22)
$ws.Range("A1").Select() | Out-Null
